$d = $word.ActiveDocument

function Replace-Text($range, $old, $new) {
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1) Hyperlinked "English" -> "Angielski" (first occurrence, paragraph 1)
Replace-Text $d.Paragraphs(1).Range "English" "Angielski"

# 2) Language list translation (avoid leading space in the search text so the
#    match does not start exactly at the hyperlink/run boundary, which would
#    otherwise cause the replaced text to incorrectly inherit the hyperlink's
#    character formatting)
Replace-Text $d.Paragraphs(1).Range "Portuguese / French / Thai / Vietnamese / Spanish" "Portugalski / Francuski / Tajski / Wietnamski / Hiszpański"

# 3) Standalone "English" heading -> "Angielski" (paragraph 3)
Replace-Text $d.Paragraphs(3).Range "English" "Angielski"

# 4) "Brief" -> "Krótki opis"
Replace-Text $d.Paragraphs(5).Range "Brief" "Krótki opis"

# 5) Brief description text
Replace-Text $d.Paragraphs(6).Range "An email sent to partners in the target country whose documents failed our verification process. It will be sent via customer.io" "Email wysłany do partnerów w docelowym kraju, których dokumenty nie przeszły naszego procesu weryfikacji. Zostanie wysłany przez customer.io"

# 6) "Target audience" -> "Docelowa grupa odbiorców"
Replace-Text $d.Paragraphs(8).Range "Target audience" "Docelowa grupa odbiorców"

# 7) Target audience description
Replace-Text $d.Paragraphs(9).Range "Invited partners who submitted wrong/incomplete documents" "Zaproszeni partnerzy, którzy przesłali błędne/niekompletne dokumenty"

# 8) "Subject line" -> "Temat"
Replace-Text $d.Paragraphs(12).Range "Subject line" "Temat"

# 9) "[EVENT NAME]" -> "[NAZWA WYDARZENIA]"
Replace-Text $d.Paragraphs(12).Range "[EVENT NAME]" "[NAZWA WYDARZENIA]"

# 10) subject suffix
Replace-Text $d.Paragraphs(12).Range " — document verification failed " " — weryfikacja dokumentu nie powiodła się "

# 11) Header
Replace-Text $d.Paragraphs(14).Range "Uh oh! Your documents couldn’t be verified" "O nie! Państwa dokumenty nie mogły zostać zweryfikowane"

# 12) "Hi " -> "Witamy "
Replace-Text $d.Paragraphs(16).Range "Hi " "Witamy "

# 13) "[PARTNER NAME]" -> "[NAZWA PARTNERA]"
Replace-Text $d.Paragraphs(16).Range "[PARTNER NAME]" "[NAZWA PARTNERA]"

# 14) Regret sentence
Replace-Text $d.Paragraphs(17).Range "We regret to inform you that your documents have failed our verification process as we found the following issues with them: " "Z przykrością informujemy, że Państwa dokumenty nie przeszły pomyślnie procesu weryfikacji, ponieważ wykryliśmy w nich następujące problemy: "

# 15) vaccination certificate bold label
Replace-Text $d.Paragraphs(18).Range "A copy of your vaccination certificate" "Kopia zaświadczenia o szczepieniu"

# 16) issue description
Replace-Text $d.Paragraphs(18).Range ": Document is unclear" ": Dokument jest nieczytelny"

# 17) [Document 2] -> [Dokument 2]
Replace-Text $d.Paragraphs(19).Range "[Document 2]" "[Dokument 2]"

# 18) Please resubmit sentence start
Replace-Text $d.Paragraphs(20).Range "Please resubmit the documents above by " "Proszę ponownie przesłać powyższe dokumenty do dnia "

# 19) sentence end
Replace-Text $d.Paragraphs(20).Range " so we can proceed with the necessary arrangements." ", abyśmy mogli dokonać niezbędnych ustaleń."

# 20) contact your country manager sentence
Replace-Text $d.Paragraphs(22).Range "If you have any questions, please contact your country manager, " "W razie jakichkolwiek pytań prosimy kontaktować się z menedżerem krajowym, "

# 21) [NAME] -> [IMIĘ]
Replace-Text $d.Paragraphs(22).Range "[NAME]" "[IMIĘ]"

# 22) ", at " -> ", pod adresem "
Replace-Text $d.Paragraphs(22).Range ", at " ", pod adresem "

# 23) [EMAIL ADDRESS] -> [ADRES EMAIL]
Replace-Text $d.Paragraphs(22).Range "[EMAIL ADDRESS]" "[ADRES EMAIL]"

# 24) " or " -> " lub " (only within this paragraph, to avoid the "live chat or WhatsApp" sentence)
Replace-Text $d.Paragraphs(22).Range " or " " lub "

# 25) [WHATSAPP NO] -> [NUMER WHATSAPP]
Replace-Text $d.Paragraphs(22).Range "[WHATSAPP NO]" "[NUMER WHATSAPP]"
